$d = $word.ActiveDocument

function New-RunsPackageXml($innerRunsXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>'
    $pkg += '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
    $pkg += '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">'
    $pkg += '<pkg:xmlData>'
    $pkg += '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">'
    $pkg += '<w:body><w:p>'
    $pkg += $innerRunsXml
    $pkg += '</w:p></w:body></w:document>'
    $pkg += '</pkg:xmlData></pkg:part></pkg:package>'
    return $pkg
}

function Replace-PlaceholderInParagraph($placeholderToken, $prefixLen, $innerRunsXml) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($placeholderToken)) {
            $target = $p
            break
        }
    }
    if ($target -eq $null) {
        throw "Paragraph containing '$placeholderToken' not found"
    }
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End
    # placeholder range = from end of the static prefix label to just before the
    # paragraph mark (End is exclusive of the paragraph mark boundary in Word's
    # Range semantics, so End-1 excludes the pilcrow).
    $placeholderRange = $d.Range($pStart + $prefixLen, $pEnd - 1)
    $xml = New-RunsPackageXml $innerRunsXml
    $placeholderRange.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1. Razão Social
# ---------------------------------------------------------------------------
$razaoRuns = @'
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:t>LEANDRO LIMA RIBEIRO FRANCA</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="36"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
'@
Replace-PlaceholderInParagraph "{{razaoSocial}}" 14 $razaoRuns

# ---------------------------------------------------------------------------
# 2. CNPJ
# ---------------------------------------------------------------------------
$cnpjRuns = @'
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:val="en-BR"/></w:rPr><w:t>60.434.015/0001-90</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="pt-BR"/></w:rPr><w:t>,</w:t></w:r>
'@
Replace-PlaceholderInParagraph "{{cnpjContratada}}" 6 $cnpjRuns

# ---------------------------------------------------------------------------
# 3. Endereço
# ---------------------------------------------------------------------------
$enderecoRuns = @'
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">RUA GOIANAZ QD 15 L 5, CONJ. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>MIRRAGE, ANAPOLIS-GO, 75070-180</w:t></w:r>
'@
Replace-PlaceholderInParagraph "{{enderecoContratada}}" 10 $enderecoRuns

Write-Output "done"
